# Two random 8-char names ("ghlktp23" and "ggu21lkq") from the unused-names
# pool on Sheet1 were consumed for a newly generated image, so:
#   1. remove those two rows from Sheet1 (the remaining names shift up)
#   2. log the two "used" events as new rows at the bottom of the "used" sheet

$wb = $excel.ActiveWorkbook

$namesSheet = $wb.Worksheets.Item("Sheet1")
$usedSheet  = $wb.Worksheets.Item("used")

# Row 1 is "ghlktp23"; deleting it shifts row 3 ("ggu21lkq") up to row 2,
# so deleting row 2 next removes "ggu21lkq" too.
$namesSheet.Rows.Item(1).Delete()
$namesSheet.Rows.Item(2).Delete()

# Append the two "used" log rows (sheet currently has data through row 21).
$usedSheet.Range("A22").Value = "ghlktp23"
$usedSheet.Range("B22").Value = "ChatGPT Image 2026年1月20日 15_06_45.png"
$usedSheet.Range("C22").Value = "2026-01-20 15:08:07"

$usedSheet.Range("A23").Value = "ggu21lkq"
$usedSheet.Range("B23").Value = "ChatGPT Image 2026年1月20日 15_06_45.png"
$usedSheet.Range("C23").Value = "2026-01-20 15:40:34"
